$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new "dbXref" column to every sheet, and populate it with hyperlinks
# to NASA science pages for the rows that have a matching external source.
# The header cell's format is copied from the existing "subclass"/"class"
# header (B1) so it reuses the same header style rather than minting a new
# (unstyled) one.
# ---------------------------------------------------------------------------

function Add-DbXrefHeader($ws, $cellRef) {
    $ws.Range($cellRef).Value = "dbXref"
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# 1) solar system -- new column E
$ws = $wb.Worksheets.Item("solar system")
Add-DbXrefHeader $ws "E1"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://science.nasa.gov/sun/", "", "", "https://science.nasa.gov/sun/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "https://science.nasa.gov/earth/", "", "", "https://science.nasa.gov/earth/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), "https://science.nasa.gov/mars/", "", "", "https://science.nasa.gov/mars/") | Out-Null

# 2) planet type -- new column D (header only)
$ws = $wb.Worksheets.Item("planet type")
Add-DbXrefHeader $ws "D1"

# 3) planetary feature -- new column E (header only)
$ws = $wb.Worksheets.Item("planetary feature")
Add-DbXrefHeader $ws "E1"

# 4) satellite -- new column E
$ws = $wb.Worksheets.Item("satellite")
Add-DbXrefHeader $ws "E1"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://science.nasa.gov/mission/hubble/", "", "", "https://science.nasa.gov/mission/hubble/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.nasa.gov/international-space-station/", "", "", "https://www.nasa.gov/international-space-station/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E7"), "https://science.nasa.gov/moon/", "", "", "https://science.nasa.gov/moon/") | Out-Null

# 5) large body object -- new column D (header only)
$ws = $wb.Worksheets.Item("large body object")
Add-DbXrefHeader $ws "D1"

# 6) small body object -- new column D
$ws = $wb.Worksheets.Item("small body object")
Add-DbXrefHeader $ws "D1"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://science.nasa.gov/solar-system/asteroids/", "", "", "https://science.nasa.gov/solar-system/asteroids/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://science.nasa.gov/solar-system/meteors-meteorites/", "", "", "https://science.nasa.gov/solar-system/meteors-meteorites/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://science.nasa.gov/solar-system/comets/", "", "", "https://science.nasa.gov/solar-system/comets/") | Out-Null

# 7) space phenomena -- new column D
$ws = $wb.Worksheets.Item("space phenomena")
Add-DbXrefHeader $ws "D1"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://science.nasa.gov/universe/black-holes/", "", "", "https://science.nasa.gov/universe/black-holes/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://science.nasa.gov/solar-system/meteors-meteorites/", "", "", "https://science.nasa.gov/solar-system/meteors-meteorites/") | Out-Null

# 8) constellation -- new column E (header only)
$ws = $wb.Worksheets.Item("constellation")
Add-DbXrefHeader $ws "E1"

Write-Output "applied dbXref columns + hyperlinks"
